# edit.ps1 -- apply the tracked content changes to Report.docx
#
# Summary of edits (see commit message: "Added Jav.doc comment for
# controller and interface classes.. and updated a class diagram"):
#   1. Merge the "We have / honoured / the principles..." declaration
#      sentence back into a single run (no text change, just run reflow).
#   2. Insert the document's "_GoBack" bookmark between "group " and
#      "them into Entity Class..." (this moves the bookmark away from
#      its old spot further down the document).
#   3. Because the bookmark moved away from "Then there will be a |
#      method which ", those two runs collapse back into a single run.
#   4. Insert " and other interface classes" after "...application
#      interface class" in the layered-architecture paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge "We have " + "honoured" + " the principles ... this work. "
#    into one run. Setting Range.Text to the *same* text is a no-op in
#    this engine, so we first swap in a throwaway placeholder to force
#    the reflow, then set the final text back.
# ---------------------------------------------------------------------
$sentence = "We have honoured the principles of academic integrity and have upheld Student Code of Academic Conduct in the completion of this work. "
$r1 = $d.Content
$found1 = $r1.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found1) {
    $r1.Text = "placeholder_merge_marker_1"
    $r1.Text = $sentence
}

# ---------------------------------------------------------------------
# 2. Split "...and group them into Entity Class..." and drop the
#    "_GoBack" bookmark in between "group " and "them into". Adding a
#    bookmark with a name that already exists elsewhere moves it here
#    (Word only keeps one bookmark per name), which also cleans up the
#    old location for step 3 below.
# ---------------------------------------------------------------------
$r2 = $d.Content
$found2 = $r2.Find.Execute("group them into", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $splitPoint = $r2.Start + 6   # just after "group "
    $bookmarkRange = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $bookmarkRange)
}

# ---------------------------------------------------------------------
# 3. Re-merge "Then there will be a " + "method which " now that the
#    bookmark that used to separate them is gone.
# ---------------------------------------------------------------------
$smsSentence = "For this extension, we can create a SMS controller class and SMS entity class. SMS controller will have a method which communicate with reservation controller to get the reservations and contact numbers for 3 days ahead compared to local date. Then there will be a method which "
$r3 = $d.Content
$found3 = $r3.Find.Execute($smsSentence, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found3) {
    $r3.Text = "placeholder_merge_marker_3"
    $r3.Text = $smsSentence
}

# ---------------------------------------------------------------------
# 4. Insert " and other interface classes" right after "...We put the
#    main application interface class" (before " in the presentation
#    layer..."). Toggling a format property on just the inserted text
#    forces it into its own run instead of merging with its neighbours.
# ---------------------------------------------------------------------
$r4 = $d.Content
$found4 = $r4.Find.Execute("We put the main application interface class", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found4) {
    $r4.Collapse(0)
    $ins = $r4.Duplicate
    $ins.InsertAfter(" and other interface classes")
    $ins.Font.Bold = $true
    $ins.Font.Bold = $false
}
